$wb = $excel.ActiveWorkbook

# =========================================================================
# Part 1: update the "总计" (summary) sheet - insert a new "2022-Q3" row
# right under the header, shifting the existing quarters down by one row.
# =========================================================================
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("A2:D2").ClearFormats()

$wsTotal.Cells.Item(2,1).Value = 0
$wsTotal.Cells.Item(2,2).Value = "2022-Q3"
$wsTotal.Cells.Item(2,3).Value = 16
$wsTotal.Cells.Item(2,4).Value = 5.3

# Copy the index-column style (s="2") from the row below onto the new A2
# cell so the whole A column keeps a uniform look.
$wsTotal.Cells.Item(3,1).Copy()
$wsTotal.Cells.Item(2,1).PasteSpecial(-4122)

# Re-number the index column (0,1,2,3,4,5) for the shifted rows.
for ($r = 3; $r -le 7; $r++) {
    $wsTotal.Cells.Item($r,1).Value = $r - 2
}

# =========================================================================
# Part 2: insert a brand-new "2022-Q3" worksheet right after "总计".
# We clone the existing "2022-Q2" sheet (so it inherits the exact same
# column layout / styles) and then overwrite its contents.
# =========================================================================
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($null, $wsTotal)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# Fund holdings for 2022-Q3 : code, name, fund size, stock position,
# position ratio, holding value (billion yuan), position rank.
$q3data = @(
    @("012930","中庚价值先锋股票","47.83","94.71","5.75","2.7502",4),
    @("920003","中金新锐股票A","17.72","89.26","5.86","1.0384",2),
    @("519033","海富通国策导向混合","15.04","91.20","2.74","0.4121",9),
    @("920923","中金新锐股票C","3.32","89.26","5.86","0.1946",2),
    @("162102","金鹰中小盘精选混合","3.17","78.28","5.12","0.1623",2),
    @("210009","金鹰核心资源混合","2.84","93.42","5.47","0.1553",5),
    @("000824","圆信永丰双红利灵活配置混合A","4.71","93.94","3.24","0.1526",10),
    @("001167","金鹰科技创新股票","2.66","94.84","5.56","0.1479",5),
    @("920002","中金精选股票A","2.95","82.28","4.03","0.1189",1),
    @("008311","圆信永丰优选价值混合A","2.29","94.30","3.53","0.0808",10),
    @("519097","新华中小市值优选混合","0.71","67.35","4.01","0.0285",5),
    @("210002","金鹰红利价值混合A","1.19","61.62","2.20","0.0262",10),
    @("000825","圆信永丰双红利灵活配置混合C","0.48","93.94","3.24","0.0156",10),
    @("016563","金鹰红利价值混合C","0.34","61.62","2.20","0.0075",10),
    @("920922","中金精选股票C","0.11","82.28","4.03","0.0044",1),
    @("008312","圆信永丰优选价值混合C","0.09","94.30","3.53","0.0032",10)
)

# The cloned sheet only has 8 rows (1 header + 7 data rows); we need
# 17 rows (1 header + 16 data rows). Extend the index column down,
# copying the style (s="2") of an existing index cell.
for ($r = 9; $r -le 17; $r++) {
    $wsQ3.Cells.Item(7,1).Copy()
    $wsQ3.Cells.Item($r,1).PasteSpecial(-4122)
}

# Make sure text-like columns (B..G) keep their text format before we
# write into them, so numeric-looking strings (e.g. "012930", "47.83")
# are not silently turned into numbers.
$wsQ3.Range("B2:G17").NumberFormat = "@"

for ($i = 0; $i -lt $q3data.Count; $i++) {
    $r = $i + 2
    $row = $q3data[$i]
    $wsQ3.Cells.Item($r,1).Value = $i
    $wsQ3.Cells.Item($r,2).Value = $row[0]
    $wsQ3.Cells.Item($r,3).Value = $row[1]
    $wsQ3.Cells.Item($r,4).Value = $row[2]
    $wsQ3.Cells.Item($r,5).Value = $row[3]
    $wsQ3.Cells.Item($r,6).Value = $row[4]
    $wsQ3.Cells.Item($r,7).Value = $row[5]
    $wsQ3.Cells.Item($r,8).Value = $row[6]
}

# Drop the temporary text-format override again so the cells end up
# without an explicit style, matching the plain data cells elsewhere.
$wsQ3.Range("B2:G17").ClearFormats()
